# Add 2022-Q4 data
#  1) "总计" sheet: insert a new data row for 2022-Q4 right after the header,
#     pushing the existing quarter rows down by one.
#  2) Insert a brand-new worksheet named "2022-Q4" right after "总计"
#     (i.e. as the new #2 sheet), holding the per-fund breakdown for the
#     quarter, in the same shape as the other quarter sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" summary sheet - shift rows 2..6 down to 3..7, then write the
#    new 2022-Q4 row into row 2.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

for ($r = 6; $r -ge 2; $r--) {
    $nr = $r + 1
    $summary.Range("B$nr").Value = $summary.Range("B$r").Value2
    $summary.Range("C$nr").Value = $summary.Range("C$r").Value2
    $summary.Range("D$nr").Value = $summary.Range("D$r").Value2
}

# Column A is just a 0-based row index (0,1,2,3,4,5) - re-stamp it (bold +
# border + centered, matching every other row in this column) now that
# there are 6 data rows instead of 5. Copy the existing A2 style (rather
# than re-deriving it from Font/Borders) so every row lands on the exact
# same style record.
for ($r = 7; $r -ge 2; $r--) {
    $summary.Range("A2").Copy($summary.Range("A$r"))
}
for ($r = 2; $r -le 7; $r++) {
    $summary.Range("A$r").Value = $r - 2
}

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 7
$summary.Range("D2").Value = 9.109999999999999

# ---------------------------------------------------------------------
# 2) New "2022-Q4" worksheet, inserted right after "总计".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

$headerRange = $q4.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# row data: A (index), B code, C name, D scale, E position, F ratio,
# G market value, H rank. B..G are stored as text in the source sheet
# (leading zeros in fund codes etc.), H is numeric.
$rows = @(
    @("320007", "诺安成长混合", "248.59", "82.42", "2.63", "6.5379", 10),
    @("002560", "诺安和鑫灵活配置混合", "31.46", "84.32", "7.30", "2.2966", 4),
    @("320022", "诺安研究精选股票", "6.42", "92.87", "2.08", "0.1335", 8),
    @("014320", "德邦半导体产业混合C", "1.52", "92.57", "4.54", "0.0690", 9),
    @("001706", "诺安积极回报灵活配置混合A", "0.52", "93.31", "7.64", "0.0397", 7),
    @("014319", "德邦半导体产业混合A", "0.37", "92.57", "4.54", "0.0168", 9),
    @("012847", "诺安积极回报灵活配置混合C", "0.18", "93.31", "7.64", "0.0138", 7)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $row = $rows[$i]
    $r = $i + 2

    $aCell = $q4.Range("A$r")
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1
    $aCell.Value = $i

    $bCell = $q4.Range("B$r")
    $bCell.NumberFormat = "@"
    $bCell.Value = $row[0]

    $q4.Range("C$r").Value = $row[1]

    $dCell = $q4.Range("D$r")
    $dCell.NumberFormat = "@"
    $dCell.Value = $row[2]

    $eCell = $q4.Range("E$r")
    $eCell.NumberFormat = "@"
    $eCell.Value = $row[3]

    $fCell = $q4.Range("F$r")
    $fCell.NumberFormat = "@"
    $fCell.Value = $row[4]

    $gCell = $q4.Range("G$r")
    $gCell.NumberFormat = "@"
    $gCell.Value = $row[5]

    $q4.Range("H$r").Value = $row[6]
}

# Keep the originally active sheet selected (adding a sheet makes it active
# by default).
$summary.Activate()
